# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of cell address -> new value for column F
$updates = @{
    "F2" = 1481
    "F3" = 959
    "F5" = 2236
    "F7" = 1368
    "F9" = 146
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
